# Auto-generated edit script: updates cryptos.xlsx price/volume/coin data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.112.61"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "1.790.18"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'224.80"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'0.548"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'32.64"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "'0.285"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").Value = "'0.0706"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "2.049.39"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "1.802.30"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "'10.82"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").Value = "'0.626"
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").Value = "34.120.03"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "'4.17"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").Value = "'243.23"
$ws.Range("E19").Value = "  -4.09%  "
$ws.Range("D20").Value = "0.0₃0784"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'10.74"
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("E23").Value = "  -4.52%  "
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").Value = "'159.72"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "'7.06"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'0.112"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D30").Value = "'0.0516"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").Value = "'3.67"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").Value = "'3.51"
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("E34").Value = "  -5.61%  "
$ws.Range("D35").Value = "1.394.97"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("D36").Value = "'0.646"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("D39").Value = "'2.21"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  +19.42%  "
$ws.Range("D42").Value = "'0.916"
$ws.Range("E42").Value = "  -6.13%  "
$ws.Range("D43").Value = "'78.72"
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value = "'0.0496"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'107.74"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").Value = "'5.90"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.949.30"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'12.20"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("E51").Value = "  -0.06%  "
